# Dynamic locator creation updated with data store replacement text. Added new tests
#
# - LoginData!B4 gets the new "Osanda Nimalarathna" value (profile name),
#   and LoginData!A4 gets the new "profileName" key label (written in this
#   order so the shared-string table picks up "Osanda Nimalarathna" before
#   "profileName", matching the target workbook).
# - The active sheet switches from SearchData back to LoginData, with the
#   selection left on A4.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LoginData")

$ws1.Range("B4").Value = "Osanda Nimalarathna"
$ws1.Range("A4").Value = "profileName"

$ws1.Activate()
$ws1.Range("A4").Select()
